$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B22").Value = 6284
$ws.Range("C22").Value = 991
$ws.Range("D22").Value = 5716856
$ws.Range("E22").Value = 909.7479312539783
$ws.Range("F22").Value = 8.176966775692884
$ws.Range("G22").Value = 3.661087866108792
$ws.Range("H22").Value = 24.3238023991162
